# pe_database.xlsx edit:
#  - Rebuild the "dcb2.0" sheet's contact-pair list with the new DCB_SH/DCB_BC
#    entries (replacing the old single "DCB2.0" placeholder row).
#  - Apply the workbook-wide default font change (Tahoma -> Calibri).
#  - Re-select / activate sheets to match the saved view state (bmw becomes
#    the active/selected tab; dcb2.0's selection moves to F16).

$wb = $excel.ActiveWorkbook

# --- 1) dcb2.0 sheet: replace the placeholder "DCB2.0" row with the new
#        contact-pair rows -------------------------------------------------
$dcb20 = $wb.Worksheets.Item("dcb2.0")

# Clear out the old single data row before rebuilding it.
$dcb20.Range("A2:B2").ClearContents()

$dcb20.Range("A2").Value = "DCB_SH to MPEC"
$dcb20.Range("B2").Value = 10
$dcb20.Range("A3").Value = "DCB_BC to MPEC"
$dcb20.Range("B3").Value = 10
# Row 5 ("DCB_BC to MPE") is entered before row 4 ("DCB_SH to MPE") so the
# shared-string table ends up with the same insertion order as the source
# workbook (BC to MPE allocated right after BC to MPEC).
$dcb20.Range("A5").Value = "DCB_BC to MPE"
$dcb20.Range("B5").Value = 10
$dcb20.Range("A4").Value = "DCB_SH to MPE"
$dcb20.Range("B4").Value = 10

$dcb20.Columns.Item(1).AutoFit()

$null = $dcb20.Activate()
$null = $dcb20.Range("F16").Select()

# --- 2) Workbook default font: Tahoma -> Calibri --------------------------
$normalStyle = $wb.Styles.Item("Normal")
$normalStyle.Font.Name = "Calibri"

# --- 3) View state: bmw becomes the active/selected sheet -----------------
$bmw = $wb.Worksheets.Item("bmw")
$null = $bmw.Activate()
$null = $bmw.Range("E13").Select()
